$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new trade row (row 8) mirroring the existing data shape (A1:I7 -> A1:I8)
# Copy row 7 formatting down to row 8 first so date/bool styles (s="1") carry over
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(8, 1).Value = 42654.743680555555
$ws.Cells.Item(8, 2).Value = $false
$ws.Cells.Item(8, 3).Value = 10079.18
$ws.Cells.Item(8, 4).Value = 10079.68
$ws.Cells.Item(8, 5).Value = 75.5
$ws.Cells.Item(8, 6).Value = 75.489998
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = -0.01
$ws.Cells.Item(8, 9).Value = $false

# Keep the best-fit column widths in sync with the newly added row
$ws.Range("A1:I8").Columns.AutoFit()
